$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.873.66"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.892.54"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.48"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4585"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("E8").Value = "  +1.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07848"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9890"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.88"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").Value = "1.848.32"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.038"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.702"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06946"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.97"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009971"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.93"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").Value = "28.872.20"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.297"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.00"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "2.091.74"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.062"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.08"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.27"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.927"
$ws.Range("E28").Value = "  +5.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.928"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.65"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09358"
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9108"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.299"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.265"
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.192"
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05766"
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02072"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.001"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.741"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5687"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1773"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.763"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.282"
$ws.Range("E44").Value = "  +7.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.88"
$ws.Range("E45").Value = "  +3.44%  "
$ws.Range("E46").Value = "  +2.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07043"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.842"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.533"
$ws.Range("E49").Value = "  +4.01%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.64"
$ws.Range("E50").Value = "  +0.81%  "
$ws.Range("E51").Value = "  -5.66%  "
